$wb = $excel.ActiveWorkbook

# Rename the existing sheet to "Attendance"
$attendance = $wb.Worksheets.Item(1)
$attendance.Name = "Attendance"

# Add a new "Key" sheet right after "Attendance"
$key = $wb.Worksheets.Add($null, $attendance)
$key.Name = "Key"

# Header row (Bunk / Name / ID) - copy from Attendance so formatting matches
$attendance.Range("A1:C1").Copy($key.Range("A1"))

# Bunk 1 staff rows - copy from Attendance (rows 6:7) so formatting (incl. quote-prefix) matches
$attendance.Range("A6:C7").Copy($key.Range("A3"))

# Bunk 2 staff rows - copy from Attendance (rows 9:10)
$attendance.Range("A9:C10").Copy($key.Range("A6"))

# Column widths to match
$key.Columns.Item(2).ColumnWidth = 83 / 6
$key.Columns.Item(3).ColumnWidth = 97 / 6

# Update selection on Attendance sheet
$attendance.Range("C28").Select() | Out-Null

# Select cell C13 on the Key sheet (leaving Key as the active/visible tab)
$key.Range("C13").Select() | Out-Null
